$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「着るのだいすき」" occupied row 555.
# Delete that entire row; Excel will automatically shift all rows below it
# (previously 556-720) up by one, so they become rows 555-719.
$ws.Rows.Item(555).Delete()
